$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 44
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "01-07-2021"
$ws.Cells.Item($row, 2).Value = 109.76
$ws.Cells.Item($row, 3).Value = 108.01
$ws.Cells.Item($row, 4).Value = 111.4
$ws.Cells.Item($row, 5).Value = 107.81
$ws.Cells.Item($row, 6).Value = 118.06
